$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B to make room for ROLL_NO
$ws.Columns.Item(2).Insert()

# Set header and values for the new ROLL_NO column
$ws.Range("B1").Value = "ROLL_NO"
$ws.Range("B2").Value = "MKY1298"
$ws.Range("B3").Value = "MKY1232"
$ws.Range("B4").Value = "MKY1222"

# Set the new column width (closest value achievable due to column-width
# quantization in the host engine; targets stored width ~21.21875)
$ws.Columns.Item(2).ColumnWidth = 20.3

# Update the active cell selection to match the new layout
$ws.Range("B4").Select()
